$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 is intentionally left blank (matches source diff which jumps from row 3 to row 5)

# Column D values are stored as text ("1"), not numbers, so force text format
# before assigning so Excel does not auto-convert them to numeric values.
$ws.Range("D5:D7").NumberFormat = "@"

$ws.Range("A5").Value = "Capital One"
$ws.Range("B5").Value = "BioCellection"
$ws.Range("C5").Value = "2020-08-31 16:28:18.472784"
$ws.Range("D5").Value = "1"

$ws.Range("A6").Value = "Heart Institute of the Caribbean"
$ws.Range("B6").Value = "Algramo-Catalyzing Reusable Packaging"
$ws.Range("C6").Value = "2020-08-31 16:29:51.732960"
$ws.Range("D6").Value = "1"

$ws.Range("A7").Value = "New Orleans Health Department"
$ws.Range("B7").Value = "BioCellection"
$ws.Range("C7").Value = "2020-08-31 16:49:09.109752"
$ws.Range("D7").Value = "1"
